$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.189.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6206"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07361"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2911"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07662"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.835.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.973"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6693"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.848"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.183.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.086.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.362"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.553"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05841"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.088"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.211"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.878"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7271"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.613"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.858"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.220.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01755"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9080"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.990.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5040"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("E48").Value = "  -4.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.148"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4028"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("E51").Value = "  +2.75%  "

